# Update the "general" sheet with two new configuration rows for
# the prior distributions of fluxes and thermodynamic quantities.
#
# The new rows are inserted above the old row 6 ("Number of exp.
# conditions ..."), pushing the existing rows (old 6-12) down to
# rows 8-14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")

# Insert two blank rows before row 6 (old rows 6-12 shift to 8-14).
$ws.Rows.Item(6).Resize(2).Insert()

# --- Row 6: Prior distribution for fluxes ---
$ws.Range("A6").Value = "Prior distribution for fluxes (uniform or normal)"
$ws.Range("B6").Value = "normal"

# --- Row 7: Prior distribution for thermodynamic quantities ---
$ws.Range("A7").Value = "Prior distribution for thermodynamic quantities (uniform or normal)"
$ws.Range("B7").Value = "normal"

# Match the formatting of the other descriptive label cells (A2:A5)
# for the new label cells in column A. Use a same-sized source range
# (A2:A3) so the paste doesn't tile past the 2-row destination.
$ws.Range("A2:A3").Copy()
$ws.Range("A6:A7").PasteSpecial(-4122)

# Format the new value cells in column B: bordered, centered,
# non-bold text (distinguishing them as selectable/editable fields).
$valueCells = $ws.Range("B6:B7")
$valueCells.Font.Name = "Calibri"
$valueCells.Font.Size = 11
$valueCells.Font.Bold = $false
$valueCells.HorizontalAlignment = -4108
$valueCells.VerticalAlignment = -4107
$valueCells.Borders.LineStyle = 1

$ws.Rows.Item(6).RowHeight = 13.8
$ws.Rows.Item(7).RowHeight = 13.8

# Make the "general" sheet the active sheet with the new rows selected,
# matching the reviewer's final view state.
$ws.Activate()
$ws.Range("A6:B7").Select()
